$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B (Coin name) updates ---
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("B45").Value = "Aave"
$ws.Range("B46").Value = "BabyDogeCoin"

# --- Column C (Link) updates ---
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"

# --- Column D (Price) updates: force text format so values like "22.00"/"0.4300" keep exact digits ---
$dCells = @("D2","D3","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D20","D21","D22","D23","D24","D25","D26","D27","D28","D32","D33","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D49","D50","D51")
foreach ($ref in $dCells) { $ws.Range($ref).NumberFormat = "@" }
$ws.Range("D2").Value = "26.414.96"
$ws.Range("D3").Value = "1.698.03"
$ws.Range("D5").Value = "218.39"
$ws.Range("D6").Value = "0.5465"
$ws.Range("D8").Value = "0.2736"
$ws.Range("D9").Value = "0.06450"
$ws.Range("D10").Value = "21.97"
$ws.Range("D11").Value = "0.07677"
$ws.Range("D12").Value = "1.695.67"
$ws.Range("D13").Value = "4.557"
$ws.Range("D14").Value = "0.5856"
$ws.Range("D15").Value = "0.000008403"
$ws.Range("D16").Value = "65.69"
$ws.Range("D17").Value = "26.484.25"
$ws.Range("D18").Value = "4.945"
$ws.Range("D20").Value = "10.98"
$ws.Range("D21").Value = "191.22"
$ws.Range("D22").Value = "6.267"
$ws.Range("D23").Value = "1.010"
$ws.Range("D24").Value = "148.69"
$ws.Range("D25").Value = "0.1311"
$ws.Range("D26").Value = "7.935"
$ws.Range("D27").Value = "15.81"
$ws.Range("D28").Value = "0.06223"
$ws.Range("D32").Value = "3.596"
$ws.Range("D33").Value = "1.686"
$ws.Range("D35").Value = "0.6161"
$ws.Range("D37").Value = "2.761"
$ws.Range("D38").Value = "0.01657"
$ws.Range("D39").Value = "1.118.09"
$ws.Range("D40").Value = "6.114"
$ws.Range("D41").Value = "0.8829"
$ws.Range("D42").Value = "1.016"
$ws.Range("D43").Value = "101.16"
$ws.Range("D44").Value = "1.849.35"
$ws.Range("D45").Value = "57.65"
$ws.Range("D46").Value = "0.00000000109"
$ws.Range("D47").Value = "8.233"
$ws.Range("D49").Value = "0.05284"
$ws.Range("D50").Value = "6.126"
$ws.Range("D51").Value = "0.4303"

# --- Column E (Volume/1h %) updates ---
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("E5").Value = "  -0.21%  "
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  +1.11%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -0.43%  "
$ws.Range("E11").Value = "  +2.32%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("E17").Value = "  +0.53%  "
$ws.Range("E18").Value = "  +0.18%  "
$ws.Range("E19").Value = "  +0.16%  "
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +0.14%  "
$ws.Range("E24").Value = "  +2.36%  "
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("E26").Value = "  +2.10%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("E29").Value = "  +2.53%  "
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("E35").Value = "  -0.84%  "
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("E37").Value = "  +1.35%  "
$ws.Range("E38").Value = "  +2.37%  "
$ws.Range("E39").Value = "  +0.88%  "
$ws.Range("E40").Value = "  -5.05%  "
$ws.Range("E41").Value = "  +0.48%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("E45").Value = "  +1.16%  "
$ws.Range("E46").Value = "  -4.95%  "
$ws.Range("E47").Value = "  +1.00%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("E51").Value = "  +0.06%  "
